{"js": "// Append a new paragraph at the end of the document body containing\n// the text \"fuck you paul\". The new paragraph naturally inherits the\n// formatting (Times New Roman) of the document's trailing paragraph\n// mark, matching the target diff.\nconst body = context.document.body;\nbody.insertParagraph(\"fuck you paul\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append a new paragraph at the very end of the document containing\n# the text \"fuck you paul\". The new paragraph naturally inherits the\n# formatting (Times New Roman) of the document's trailing paragraph\n# mark, matching the target diff.\n$d = $word.ActiveDocument\n$range = $d.Content\n$range.Collapse(0)  # wdCollapseEnd\n$range.InsertParagraphAfter()\n$range.Collapse(0)  # wdCollapseEnd\n$range.Text = \"fuck you paul\"\n"}
